$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$pm = [string][char]0x00B1

$ws.Range("B2").Value = "0.018" + $pm + "0.001"
$ws.Range("C2").Value = "0.213" + $pm + "0.003"

$ws.Range("B3").Value = "0.063" + $pm + "0.004"
$ws.Range("C3").Value = "0.232" + $pm + "0.014"

$ws.Range("B4").Value = "0.872" + $pm + "0.003"
$ws.Range("C4").Value = "0.494" + $pm + "0.043"

$ws.Range("B5").Value = "0.993" + $pm + "0.000"
$ws.Range("C5").Value = "0.436" + $pm + "0.038"

$ws.Range("B6").Value = "0.979" + $pm + "0.006"
$ws.Range("C6").Value = "0.794" + $pm + "0.020"

$ws.Range("B7").Value = "0.949" + $pm + "0.003"
$ws.Range("C7").Value = "0.339" + $pm + "0.030"

$ws.Range("B8").Value = "0.009" + $pm + "0.000"
$ws.Range("C8").Value = "0.182" + $pm + "0.008"

$ws.Range("B9").Value = "0.162" + $pm + "0.030"
$ws.Range("C9").Value = "0.226" + $pm + "0.019"

$ws.Range("B10").Value = "0.664" + $pm + "0.040"
$ws.Range("C10").Value = "0.391" + $pm + "0.030"
